# Weekly refresh: insert this week's new record at the top of the data
# block (row 42) and push the existing history down by one row, matching
# the "Fruta / hortaliza, semanal" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 42..68 down to 43..69, creating a blank row 42.
$ws.Rows.Item(42).Insert()

# Populate the new weekly record in row 42.
$ws.Range("A42").Value = 9
$ws.Range("B42").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C42").Value = "Metropolitana"
$ws.Range("D42").Value = 44606
$ws.Range("E42").Value = 13
$ws.Range("F42").Value = 100114002
$ws.Range("G42").Value = "Camote"
$ws.Range("H42").Value = "Sin especificar"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 970
$ws.Range("K42").Value = 10000
$ws.Range("L42").Value = 11000
$ws.Range("M42").Value = 10495
$ws.Range("N42").Value = '$/malla 18 kilos'
$ws.Range("O42").Value = "Perú"
$ws.Range("P42").Value = 583
$ws.Range("Q42").Value = 18
$ws.Range("R42").Value = "Hortaliza"
